# Updates the cryptos list (Price / Volume(1h) columns) to the latest
# scraped figures, and fixes two rows whose coin data had been swapped
# (RenderToken/ImmutableX and ThetaToken/ApeXProtocol).
#
# Note: several "Price" values look like plain numbers (e.g. "573.11").
# Setting such a string via .Value would make Excel auto-convert the
# cell to a numeric type (losing the original text formatting and
# introducing floating point artifacts). To keep these cells as text -
# matching the workbook's original inlineStr/shared-string representation -
# we temporarily force a Text number format before assigning the value,
# then reset the style back to "Normal" so no stray style id is left
# on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.206.10"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "3.516.16"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").Value = "3.505.48"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.658"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "4.079.49"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("E16").Value = "  -2.73%  "
$ws.Range("D17").Value = "3.520.29"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "69.140.84"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "546.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +14.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("E31").Value = "  -6.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "558.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("E37").Value = "  +7.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.402"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "0.0₃0766"
$ws.Range("E40").Value = "  -5.03%  "
$ws.Range("E41").Value = "  -4.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("D44").Value = "3.271.26"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.97%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0445"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("E49").Value = "  -4.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.17%  "
